# Apply the schedule.xlsx update: two new log entries (F4 "work on navi
# bar" and C7 "/") plus tidy A11's date style so it matches the other
# date cells in column A, then leave the selection where the user's
# cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note tucked into column F on the "Added music" row.
$ws.Range("F4").Value = "work on navi bar"

# Row 7 gains the same "/" separator cell used by the other entries,
# center-aligned like its neighbours in column C.
$ws.Range("C7").Value = "/"
$ws.Range("C7").HorizontalAlignment = -4108

# A11's date was carrying its own one-off style; re-apply the shared
# YYYY-MM-DD format so it collapses onto the same style as the rest of
# column A.
$ws.Range("A11").NumberFormat = "YYYY\-MM\-DD;@"

# Leave the selection where editing left off.
$ws.Range("G12").Select() | Out-Null
